$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 992.6111
$ws.Range("I18").Value = 933.2222
$ws.Range("J18").Value = 1052
$ws.Range("K18").Value = 933.2222
$ws.Range("L18").Value = 1052
$ws.Range("M18").Value = -649.2222
$ws.Range("N18").Value = -1620
$ws.Range("H88").Value = 2395
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 2395
$ws.Range("K88").Value = 0
$ws.Range("M88").Value = 2395
$ws.Range("N88").Value = -3207
$ws.Range("H91").Value = 2395
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 2395
$ws.Range("K91").Value = 0
$ws.Range("M91").Value = 2395
$ws.Range("N91").Value = -5203
$ws.Range("H111").Value = 5944
$ws.Range("I111").Value = 6451.2144
$ws.Range("K111").Value = 19353.6432
$ws.Range("M111").Value = -16286.6432
$ws.Range("H125").Value = 24444
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 24444
$ws.Range("K125").Value = 0
$ws.Range("M125").Value = 219996
$ws.Range("N125").Value = -224916
$ws.Range("H137").Value = 3941.25
$ws.Range("I137").Value = 1537.25
$ws.Range("K137").Value = 4611.75
$ws.Range("M137").Value = -2061.75
$ws.Range("H138").Value = 2902.75
$ws.Range("J138").Value = 5958.5
$ws.Range("L138").Value = 17875.5
$ws.Range("N138").Value = -28155.5
$ws.Range("L88").ClearContents()
$ws.Range("L91").ClearContents()
$ws.Range("L125").ClearContents()

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3539.8
$ws.Range("I2").Value = 2849.5
$ws.Range("K2").Value = 2849.5
$ws.Range("M2").Value = -2736.5
$ws.Range("H45").Value = 4254.1113
$ws.Range("I45").Value = 5572.8335
$ws.Range("K45").Value = 5572.8335
$ws.Range("M45").Value = -5195.8335
$ws.Range("H61").Value = 1899.5
$ws.Range("I61").Value = 1399
$ws.Range("K61").Value = 1399
$ws.Range("M61").Value = -1187
$ws.Range("H74").Value = 2155.4
$ws.Range("I74").Value = 2899
$ws.Range("J74").Value = 1659.6666
$ws.Range("K74").Value = 2899
$ws.Range("L74").Value = 1659.6666
$ws.Range("M74").Value = -2025
$ws.Range("N74").Value = -3407.6666
$ws.Range("H77").Value = 2155.4
$ws.Range("I77").Value = 2899
$ws.Range("J77").Value = 1659.6666
$ws.Range("K77").Value = 14495
$ws.Range("L77").Value = 8298.333000000001
$ws.Range("M77").Value = -10127
$ws.Range("N77").Value = -17034.333
$ws.Range("H116").Value = 3539.8
$ws.Range("I116").Value = 2849.5
$ws.Range("K116").Value = 2849.5
$ws.Range("M116").Value = -555.5
$ws.Range("H132").Value = 983.25
$ws.Range("I132").Value = 911.5
$ws.Range("J132").Value = 1055
$ws.Range("K132").Value = 2734.5
$ws.Range("L132").Value = 3165
$ws.Range("M132").Value = -204.5
$ws.Range("N132").Value = -8225
$ws.Range("H136").Value = 1899.5
$ws.Range("I136").Value = 1399
$ws.Range("K136").Value = 4197
$ws.Range("M136").Value = -1647

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3539.8
$ws.Range("I3").Value = 2849.5
$ws.Range("K3").Value = 2849.5
$ws.Range("M3").Value = -2735.5
$ws.Range("H132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("N132").Value = 0
$ws.Range("H134").Value = 0
$ws.Range("I134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("L132").ClearContents()
$ws.Range("M134").ClearContents()

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2908.4167
$ws.Range("I16").Value = 3400.5
$ws.Range("J16").Value = 2416.3333
$ws.Range("K16").Value = 3400.5
$ws.Range("L16").Value = 2416.3333
$ws.Range("M16").Value = -3113.5
$ws.Range("N16").Value = -2990.3333
$ws.Range("H18").Value = 15000
$ws.Range("J18").Value = 15000
$ws.Range("L18").Value = 15000
$ws.Range("N18").Value = -15460
$ws.Range("H94").Value = 106050.27
$ws.Range("I94").Value = 189776.83
$ws.Range("K94").Value = 189776.83
$ws.Range("M94").Value = -189325.83
$ws.Range("H113").Value = 2908.4167
$ws.Range("I113").Value = 3400.5
$ws.Range("J113").Value = 2416.3333
$ws.Range("K113").Value = 3400.5
$ws.Range("L113").Value = 2416.3333
$ws.Range("M113").Value = -1230.5
$ws.Range("N113").Value = -6756.3333
$ws.Range("H115").Value = 43999.5
$ws.Range("J115").Value = 43999.5
$ws.Range("L115").Value = 43999.5
$ws.Range("N115").Value = -46349.5

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 173
$ws.Range("J92").Value = 209.5
$ws.Range("L92").Value = 628.5
$ws.Range("N92").Value = -3124.5
$ws.Range("H95").Value = 8883
$ws.Range("J95").Value = 8883
$ws.Range("L95").Value = 26649
$ws.Range("N95").Value = -30767
$ws.Range("H97").Value = 550
$ws.Range("I97").Value = 0
$ws.Range("K97").Value = 0
$ws.Range("H98").Value = 3058.5
$ws.Range("I98").Value = 2754.6667
$ws.Range("J98").Value = 3970
$ws.Range("K98").Value = 8264.000100000001
$ws.Range("L98").Value = 11910
$ws.Range("M98").Value = -6766.000100000001
$ws.Range("N98").Value = -14906
$ws.Range("M97").ClearContents()

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 6385.0415
$ws.Range("I132").Value = 6488.7393
$ws.Range("K132").Value = 19466.2179
$ws.Range("M132").Value = -16936.2179

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H10").Value = 333999.66
$ws.Range("I10").Value = 333999.66
$ws.Range("K10").Value = 333999.66
$ws.Range("M10").Value = -333859.66
$ws.Range("H46").Value = 1618.8462
$ws.Range("I46").Value = 1929
$ws.Range("J46").Value = 1525.8
$ws.Range("K46").Value = 1929
$ws.Range("L46").Value = 1525.8
$ws.Range("M46").Value = -1741
$ws.Range("N46").Value = -1901.8
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("H82").Value = 2719.2
$ws.Range("J82").Value = 2719.2
$ws.Range("L82").Value = 2719.2
$ws.Range("N82").Value = -3441.2
$ws.Range("H85").Value = 2719.2
$ws.Range("J85").Value = 2719.2
$ws.Range("L85").Value = 2719.2
$ws.Range("N85").Value = -5215.2
$ws.Range("H100").Value = 4766.25
$ws.Range("I100").Value = 5121.8887
$ws.Range("K100").Value = 5121.8887
$ws.Range("M100").Value = -4580.8887
$ws.Range("H136").Value = 2199.5
$ws.Range("I136").Value = 1400
$ws.Range("J136").Value = 2999
$ws.Range("K136").Value = 4200
$ws.Range("L136").Value = 8997
$ws.Range("M136").Value = -1650
$ws.Range("N136").Value = -14097
$ws.Range("M63").ClearContents()
$ws.Range("M66").ClearContents()

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H10").Value = 900
$ws.Range("I10").Value = 900
$ws.Range("K10").Value = 900
$ws.Range("M10").Value = -731
$ws.Range("H70").Value = 34999
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("H73").Value = 34999
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("H136").Value = 915.6923
$ws.Range("I136").Value = 908.4286
$ws.Range("J136").Value = 924.1667
$ws.Range("K136").Value = 2725.2858
$ws.Range("L136").Value = 2772.5001
$ws.Range("M136").Value = -175.2857999999997
$ws.Range("N136").Value = -7872.5001
$ws.Range("N70").ClearContents()
$ws.Range("N73").ClearContents()
